$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.773.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.144.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.142.33"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.524"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.87%  "
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000257"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.657.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.901.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.140.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.55%  "
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "500.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.709"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.78"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.53"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("E37").Value = "  +3.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "470.73"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0413"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.015.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.116"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.281"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0573"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.114"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "117.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.94%  "
